# Weekly crypto-tracker refresh (GitHub Actions scraper run).
# For rows whose new value is a bare number/digit string, Excel would
# auto-convert it to a numeric cell on assignment; the source sheet keeps
# every data column (B:G) as text, so Set-TextValue forces it back to text
# the same way typing a leading apostrophe in Excel's UI would.
function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
Set-TextValue $ws.Range("D2") "245.45"
Set-TextValue $ws.Range("G2") "2"

# row 3
Set-TextValue $ws.Range("D3") "22.12"
Set-TextValue $ws.Range("G3") "2"

# row 4
Set-TextValue $ws.Range("D4") "5.360"
Set-TextValue $ws.Range("G4") "2"

# row 5
Set-TextValue $ws.Range("D5") "0.05899"
Set-TextValue $ws.Range("G5") "2"

# row 6
Set-TextValue $ws.Range("G6") "2"

# row 7
Set-TextValue $ws.Range("D7") "6.387"
Set-TextValue $ws.Range("G7") "2"

# row 8
Set-TextValue $ws.Range("D8") "0.8092"
Set-TextValue $ws.Range("G8") "2"

# row 9
Set-TextValue $ws.Range("D9") "0.9624"
Set-TextValue $ws.Range("G9") "2"

# row 10 -> WazirX
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D10") "0.1427"
$ws.Range("E10").Value = "9WazirXWRX"
Set-TextValue $ws.Range("G10") "2"

# row 11 -> MandalaExchangeToken
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D11") "0.07383"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
Set-TextValue $ws.Range("G11") "2"

# row 12 -> LiechtensteinCryptoassetsExchange
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D12") "0.03448"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue $ws.Range("G12") "2"

# row 13 -> BitrueCoin
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.03035"
$ws.Range("E13").Value = "12BitrueCoinBTR"
Set-TextValue $ws.Range("G13") "2"

# row 14 -> MCDex
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D14") "4.425"
$ws.Range("E14").Value = "13MCDexMCB"
Set-TextValue $ws.Range("G14") "2"

# row 15 -> BitMartToken
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D15") "0.09401"
$ws.Range("E15").Value = "14BitMartTokenBMX"
Set-TextValue $ws.Range("G15") "2"

# row 16 -> BitForexToken
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D16") "0.001587"
$ws.Range("E16").Value = "15BitForexTokenBF"
Set-TextValue $ws.Range("G16") "2"

# row 17 -> CoinExToken
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D17") "0.04838"
$ws.Range("E17").Value = "16CoinExTokenCET"
Set-TextValue $ws.Range("G17") "2"

# row 18 -> One
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D18") "0.0005902"
$ws.Range("E18").Value = "17OneONEWorstin24h"
Set-TextValue $ws.Range("G18") "2"

# row 19
Set-TextValue $ws.Range("D19") "0.006180"
Set-TextValue $ws.Range("G19") "2"

# row 20
Set-TextValue $ws.Range("D20") "0.004074"
Set-TextValue $ws.Range("G20") "2"

# row 21
Set-TextValue $ws.Range("D21") "0.0009835"
Set-TextValue $ws.Range("G21") "2"

# row 22
Set-TextValue $ws.Range("D22") "0.00009703"
Set-TextValue $ws.Range("G22") "2"

# row 23
Set-TextValue $ws.Range("D23") "3.704"
Set-TextValue $ws.Range("G23") "2"

# row 24
Set-TextValue $ws.Range("D24") "2.188"
Set-TextValue $ws.Range("G24") "2"

# row 25
Set-TextValue $ws.Range("G25") "2"

# row 26
Set-TextValue $ws.Range("D26") "0.1340"
Set-TextValue $ws.Range("G26") "2"

# row 27
Set-TextValue $ws.Range("G27") "2"

# row 28
Set-TextValue $ws.Range("G28") "2"

# row 29
Set-TextValue $ws.Range("G29") "2"

# row 30
Set-TextValue $ws.Range("G30") "2"

# row 31
Set-TextValue $ws.Range("G31") "2"

# row 32
Set-TextValue $ws.Range("G32") "2"

# row 33
Set-TextValue $ws.Range("G33") "2"

# row 34
Set-TextValue $ws.Range("G34") "2"

# row 35
Set-TextValue $ws.Range("G35") "2"

# row 36
Set-TextValue $ws.Range("G36") "2"

# row 37
Set-TextValue $ws.Range("G37") "2"

# row 38
Set-TextValue $ws.Range("G38") "2"

# row 39
Set-TextValue $ws.Range("G39") "2"

# row 40
Set-TextValue $ws.Range("D40") "0.03925"
Set-TextValue $ws.Range("G40") "2"

# row 41
Set-TextValue $ws.Range("D41") "0.006622"
Set-TextValue $ws.Range("G41") "2"

# row 42
Set-TextValue $ws.Range("D42") "0.1072"
Set-TextValue $ws.Range("G42") "2"

# row 43
Set-TextValue $ws.Range("D43") "0.003001"
Set-TextValue $ws.Range("G43") "2"

# row 44
Set-TextValue $ws.Range("D44") "0.005778"
Set-TextValue $ws.Range("G44") "2"

# row 45
Set-TextValue $ws.Range("D45") "0.00005300"
Set-TextValue $ws.Range("G45") "2"

# row 46
Set-TextValue $ws.Range("G46") "2"

# row 47
Set-TextValue $ws.Range("D47") "0.6902"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
Set-TextValue $ws.Range("G47") "2"

# row 48
Set-TextValue $ws.Range("D48") "0.05085"
$ws.Range("E48").Value = "47BOLOBOLO"
Set-TextValue $ws.Range("G48") "2"

# row 49
Set-TextValue $ws.Range("G49") "2"

# row 50
Set-TextValue $ws.Range("D50") "0.01010"
Set-TextValue $ws.Range("G50") "2"

# row 51
Set-TextValue $ws.Range("G51") "2"
